# Aktualisierung des Zeitplans fuer die finale Praesentation (vh)
#
# Marks a batch of tasks in the "Aufgabenliste Projekt 1" sheet as
# completed (% erledigt = 100%) and records their actual completion
# dates ("tatsaechliche Fertigstellung"). Also fixes a stray "%" typed
# into the final milestone row, and moves the sheet's active-cell
# selection forward to E50.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Aufgabenliste Projekt 1")

# --- Mark tasks as 100% done -------------------------------------------
# (the "Fortschritt"/F column is a calculated table column that re-derives
# itself from "% erledigt"/E automatically once E is written)
$doneRows = 39,41,42,43,44,45,46,47,48,49
foreach ($r in $doneRows) {
    $ws.Range("E$r").Value = 1
}

# --- Record actual completion dates for the newly finished tasks -------
# Copy number formatting from a sibling cell that already carries the
# "tatsaechliche Fertigstellung" date style (xlPasteFormats = -4122),
# then fill in the date value (serials, same as Excel stores internally).
$xlPasteFormats = -4122
$excel.CutCopyMode = $false
$ws.Range("H39").Copy()
foreach ($cellRef in @("H42","H43","H44","H46","H47","H48","H49")) {
    $ws.Range($cellRef).PasteSpecial($xlPasteFormats)
}
$excel.CutCopyMode = $false
$ws.Range("H23").Copy()
$ws.Range("H45").PasteSpecial($xlPasteFormats)
$excel.CutCopyMode = $false

$completionDates = @{
    "H42" = 42345  # 07.12.2015
    "H43" = 42345  # 07.12.2015
    "H44" = 42345  # 07.12.2015
    "H45" = 42344  # 06.12.2015
    "H46" = 42346  # 08.12.2015
    "H47" = 42345  # 07.12.2015
    "H48" = 42346  # 08.12.2015
    "H49" = 42347  # 09.12.2015
}
foreach ($cellRef in $completionDates.Keys) {
    $ws.Range($cellRef).Value = $completionDates[$cellRef]
}

# --- Fix the stray "%" typed into the final milestone row --------------
$ws.Range("E52").Value = "%"

# --- Move the active selection ------------------------------------------
$ws.Range("E50").Select()
